$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC33_price_Verification_on_Cart")

# Delete row 23 (VERIFY_TEXT_PRESENT / OrderSummary row) - this row is removed entirely,
# causing the rows below it to shift up by one.
$ws.Rows("23:23").Delete()

# Delete what are now rows 28:32 (originally rows 29:33, the Price verification rows),
# which are removed entirely so the data now ends at row 27.
$ws.Rows("28:32").Delete()

# Update the active selection on the sheet to A2 (no more scrolled/selected F20 state).
$ws.Activate()
$ws.Range("A2").Select()
